# Insert a new data row before row 320 (shifts existing rows 320:385 down to 321:386)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new record's data.
$ws.Range("A320").Value = 4
$ws.Range("B320").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C320").Value = "Los Lagos"
$ws.Range("D320").Value = 44995
$ws.Range("E320").Value = 10
$ws.Range("F320").Value = 100112021
$ws.Range("G320").Value = "Ají"
$ws.Range("H320").Value = "Inferno"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 180
$ws.Range("K320").Value = 18000
$ws.Range("L320").Value = 18000
$ws.Range("M320").Value = 18000
$ws.Range("N320").Value = "$/caja 10 kilos"
$ws.Range("O320").Value = "Región de Arica y Parinacota"
$ws.Range("P320").Value = 1800
$ws.Range("Q320").Value = 10
$ws.Range("R320").Value = "Hortaliza"
